$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.912.21'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '3.309.07'
$ws.Range('E3').Value = '  +5.87%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.21'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.43%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.308.38'
$ws.Range('E8').Value = '  +6.09%  '
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('E11').Value = '  +4.86%  '
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.93'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').Value = '3.850.95'
$ws.Range('E15').Value = '  +5.88%  '
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '3.306.52'
$ws.Range('E17').Value = '  +5.95%  '
$ws.Range('D18').Value = '64.011.59'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.90'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '482.39'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.32'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('E22').Value = '  +6.26%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.03'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.54'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +3.90%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.22'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.29'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.29'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.66%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +3.79%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.40'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.90%  '
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.56'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('E36').Value = '  +2.75%  '
$ws.Range('D37').Value = '0.0₃0766'
$ws.Range('E37').Value = '  +7.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '53.34'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.64%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0400'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.05%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '430.78'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.91%  '
$ws.Range('E41').Value = '  +5.62%  '
$ws.Range('D42').Value = '3.034.22'
$ws.Range('E42').Value = '  +5.09%  '
$ws.Range('E43').Value = '  +1.88%  '
$ws.Range('E44').Value = '  -6.20%  '
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.25'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.13%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.52'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.54%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.115'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '35.53'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +13.47%  '
